$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "42.442.56"
$ws.Range("E2").Value = "  -0.72%  "
$ws.Range("D3").Value = "2.253.08"
$ws.Range("E3").Value = "  -0.56%  "
$ws.Range("E4").Value = "  +0.05%  "
Set-TextValue $ws "D5" "246.72"
$ws.Range("E5").Value = "  -1.13%  "
$ws.Range("E6").Value = "  -0.85%  "
Set-TextValue $ws "D7" "76.00"
$ws.Range("E7").Value = "  -2.93%  "
$ws.Range("E8").Value = "  +0.09%  "
Set-TextValue $ws "D9" "0.622"
$ws.Range("E9").Value = "  -3.79%  "
Set-TextValue $ws "D10" "43.88"
$ws.Range("E10").Value = "  +8.48%  "
Set-TextValue $ws "D11" "0.0952"
$ws.Range("E11").Value = "  -1.02%  "
Set-TextValue $ws "D12" "7.27"
$ws.Range("E12").Value = "  -0.79%  "
Set-TextValue $ws "D13" "0.104"
$ws.Range("E13").Value = "  -1.56%  "
$ws.Range("D14").Value = "2.588.08"
$ws.Range("E14").Value = "  -0.62%  "
Set-TextValue $ws "D15" "14.65"
$ws.Range("E15").Value = "  -2.07%  "
Set-TextValue $ws "D16" "0.858"
$ws.Range("E16").Value = "  -0.48%  "
$ws.Range("D17").Value = "2.235.76"
$ws.Range("E17").Value = "  -1.37%  "
$ws.Range("D18").Value = "42.229.05"
$ws.Range("E18").Value = "  -0.84%  "
Set-TextValue $ws "D19" "0.0000102"
$ws.Range("E19").Value = "  +3.61%  "
Set-TextValue $ws "D20" "6.19"
$ws.Range("E20").Value = "  +0.11%  "
Set-TextValue $ws "D21" "72.35"
$ws.Range("E21").Value = "  +0.79%  "
Set-TextValue $ws "D22" "2.23"
$ws.Range("E22").Value = "  +3.55%  "
Set-TextValue $ws "D23" "232.23"
$ws.Range("E23").Value = "  -0.31%  "
Set-TextValue $ws "D24" "9.26"
$ws.Range("E24").Value = "  +37.90%  "
$ws.Range("E25").Value = "  +0.10%  "
Set-TextValue $ws "D26" "11.49"
$ws.Range("E26").Value = "  +1.36%  "
$ws.Range("E27").Value = "  -4.75%  "
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("E29").Value = "  +4.17%  "
Set-TextValue $ws "D30" "168.95"
$ws.Range("E30").Value = "  -0.61%  "
Set-TextValue $ws "D31" "20.73"
$ws.Range("E31").Value = "  -0.03%  "
Set-TextValue $ws "D32" "0.0829"
$ws.Range("E32").Value = "  -3.11%  "
$ws.Range("E33").Value = "  +0.39%  "
Set-TextValue $ws "D34" "30.84"
$ws.Range("E34").Value = "  +1.93%  "
Set-TextValue $ws "D35" "5.40"
$ws.Range("E35").Value = "  +14.41%  "
Set-TextValue $ws "D36" "0.126"
$ws.Range("E36").Value = "  -0.58%  "
$ws.Range("E37").Value = "  -1.33%  "
$ws.Range("E38").Value = "  +5.31%  "
Set-TextValue $ws "D39" "13.83"
$ws.Range("E39").Value = "  +4.22%  "
$ws.Range("E40").Value = "  -1.93%  "
Set-TextValue $ws "D41" "5.82"
$ws.Range("E41").Value = "  -2.81%  "
Set-TextValue $ws "D42" "63.40"
$ws.Range("E42").Value = "  +3.68%  "
$ws.Range("E43").Value = "  -0.98%  "
Set-TextValue $ws "D44" "109.23"
$ws.Range("E44").Value = "  -4.44%  "
Set-TextValue $ws "D45" "8.80"
$ws.Range("E45").Value = "  -1.01%  "
$ws.Range("E46").Value = "  +1.62%  "
$ws.Range("E47").Value = "  -0.14%  "
$ws.Range("E48").Value = "  +2.48%  "
$ws.Range("E49").Value = "  -1.61%  "
Set-TextValue $ws "D50" "2.33"
$ws.Range("E50").Value = "  +3.96%  "
Set-TextValue $ws "D51" "4.17"
$ws.Range("E51").Value = "  -8.05%  "
